# Sprint 2 Backlog update
# - Update Status column (B) values to reflect the latest sprint progress
# - Update the active selection to C29 (as left by the author on save)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value  = "Done"
$ws.Range("B3").Value  = "Done"
$ws.Range("B4").Value  = "Done"
$ws.Range("B9").Value  = "Reviewing"
$ws.Range("B10").Value = "Reviewing"
$ws.Range("B12").Value = "Done"
$ws.Range("B13").Value = "Done"
$ws.Range("B14").Value = "Done"
$ws.Range("B16").Value = "Done"

$ws.Range("C29").Select()
